# Apply stackup updates recommended by fab engineer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8-layer")

# Update trace-width / calculated-Z0 values on rows 2, 4, 14, 16
$ws.Cells.Item(2, 8).Value = 5
$ws.Cells.Item(2, 9).Value = 88

$ws.Cells.Item(4, 8).Value = 5
$ws.Cells.Item(4, 9).Value = 49.9

$ws.Cells.Item(14, 8).Value = 5
$ws.Cells.Item(14, 9).Value = 49.9

$ws.Cells.Item(16, 8).Value = 5
$ws.Cells.Item(16, 9).Value = 88

# Narrow column E (Er) to match updated layout
$ws.Columns.Item(5).ColumnWidth = 6.451822916666667

# Update the active selection to reflect where the engineer left off
[void]$ws.Range("G19").Select()

# Remove the now-unused empty "Sheet3"
$excel.DisplayAlerts = $false
$sheet3 = $wb.Worksheets.Item("Sheet3")
[void]$sheet3.Delete()
$excel.DisplayAlerts = $true
